# Correção nos dados: as linhas 5 ("situação do domicílio") e, após a
# primeira remoção, a linha 7 ("grandes regiões e unidades da federação")
# eram apenas rótulos de cabeçalho sem dados numéricos próprios - os
# valores de urbana/rural/norte/... estavam, por engano, uma linha abaixo
# de onde deveriam estar. Removendo essas duas linhas inteiras, os dados
# numéricos (colunas B:D) sobem e ficam alinhados com o rótulo correto
# em A, e as duas últimas linhas (antigas 38 e 39) deixam de existir.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove a linha do cabeçalho "situação do domicílio" - os dados de
# "urbana", "rural", etc. sobem uma posição.
$ws.Rows(5).Delete() | Out-Null

# Depois da remoção acima, o antigo cabeçalho "grandes regiões e
# unidades da federação" (linha 8 original) passou a ser a linha 7.
# Remove essa linha também, para que os dados de "norte", "rondônia",
# etc. subam e fiquem alinhados com seus rótulos.
$ws.Rows(7).Delete() | Out-Null
